$wb = $excel.ActiveWorkbook

# 1. Insert a new worksheet "2022-Q1" before the "总计" (total) sheet.
#    Copying an existing quarter sheet preserves the column formatting/styles,
#    then we overwrite its contents with the 2022-Q1 figures.
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$srcSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# NOTE: re-resolve the "总计" sheet by name - after Copy/rename the old
# $totalSheet variable can end up pointing at the newly inserted sheet instead.
$totalSheet = $wb.Worksheets.Item("总计")

# Remove any left-over rows copied from the source sheet beyond our new data range.
$newSheet.Range("A35:H42").Clear()

# Fill header row (already correct from the copied sheet) and data rows.
# Row 2
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).NumberFormat = "@"
$newSheet.Cells.Item(2, 2).Value = "001645"
$newSheet.Cells.Item(2, 2).Style = "Normal"
$newSheet.Cells.Item(2, 3).NumberFormat = "@"
$newSheet.Cells.Item(2, 3).Value = "国泰大健康股票A"
$newSheet.Cells.Item(2, 3).Style = "Normal"
$newSheet.Cells.Item(2, 4).NumberFormat = "@"
$newSheet.Cells.Item(2, 4).Value = "34.81"
$newSheet.Cells.Item(2, 4).Style = "Normal"
$newSheet.Cells.Item(2, 5).NumberFormat = "@"
$newSheet.Cells.Item(2, 5).Value = "90.83"
$newSheet.Cells.Item(2, 5).Style = "Normal"
$newSheet.Cells.Item(2, 6).NumberFormat = "@"
$newSheet.Cells.Item(2, 6).Value = "7.96"
$newSheet.Cells.Item(2, 6).Style = "Normal"
$newSheet.Cells.Item(2, 7).NumberFormat = "@"
$newSheet.Cells.Item(2, 7).Value = "2.7709"
$newSheet.Cells.Item(2, 7).Style = "Normal"
$newSheet.Cells.Item(2, 8).Value = 3

# Row 3
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).NumberFormat = "@"
$newSheet.Cells.Item(3, 2).Value = "519995"
$newSheet.Cells.Item(3, 2).Style = "Normal"
$newSheet.Cells.Item(3, 3).NumberFormat = "@"
$newSheet.Cells.Item(3, 3).Value = "长信金利趋势混合"
$newSheet.Cells.Item(3, 3).Style = "Normal"
$newSheet.Cells.Item(3, 4).NumberFormat = "@"
$newSheet.Cells.Item(3, 4).Value = "56.80"
$newSheet.Cells.Item(3, 4).Style = "Normal"
$newSheet.Cells.Item(3, 5).NumberFormat = "@"
$newSheet.Cells.Item(3, 5).Value = "86.00"
$newSheet.Cells.Item(3, 5).Style = "Normal"
$newSheet.Cells.Item(3, 6).NumberFormat = "@"
$newSheet.Cells.Item(3, 6).Value = "3.84"
$newSheet.Cells.Item(3, 6).Style = "Normal"
$newSheet.Cells.Item(3, 7).NumberFormat = "@"
$newSheet.Cells.Item(3, 7).Value = "2.1811"
$newSheet.Cells.Item(3, 7).Style = "Normal"
$newSheet.Cells.Item(3, 8).Value = 3

# Row 4
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).NumberFormat = "@"
$newSheet.Cells.Item(4, 2).Value = "008065"
$newSheet.Cells.Item(4, 2).Style = "Normal"
$newSheet.Cells.Item(4, 3).NumberFormat = "@"
$newSheet.Cells.Item(4, 3).Value = "汇添富中盘积极成长混合A"
$newSheet.Cells.Item(4, 3).Style = "Normal"
$newSheet.Cells.Item(4, 4).NumberFormat = "@"
$newSheet.Cells.Item(4, 4).Value = "41.22"
$newSheet.Cells.Item(4, 4).Style = "Normal"
$newSheet.Cells.Item(4, 5).NumberFormat = "@"
$newSheet.Cells.Item(4, 5).Value = "73.48"
$newSheet.Cells.Item(4, 5).Style = "Normal"
$newSheet.Cells.Item(4, 6).NumberFormat = "@"
$newSheet.Cells.Item(4, 6).Value = "3.68"
$newSheet.Cells.Item(4, 6).Style = "Normal"
$newSheet.Cells.Item(4, 7).NumberFormat = "@"
$newSheet.Cells.Item(4, 7).Value = "1.5169"
$newSheet.Cells.Item(4, 7).Style = "Normal"
$newSheet.Cells.Item(4, 8).Value = 4

# Row 5
$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).NumberFormat = "@"
$newSheet.Cells.Item(5, 2).Value = "020001"
$newSheet.Cells.Item(5, 2).Style = "Normal"
$newSheet.Cells.Item(5, 3).NumberFormat = "@"
$newSheet.Cells.Item(5, 3).Value = "国泰金鹰增长灵活配置混合"
$newSheet.Cells.Item(5, 3).Style = "Normal"
$newSheet.Cells.Item(5, 4).NumberFormat = "@"
$newSheet.Cells.Item(5, 4).Value = "17.38"
$newSheet.Cells.Item(5, 4).Style = "Normal"
$newSheet.Cells.Item(5, 5).NumberFormat = "@"
$newSheet.Cells.Item(5, 5).Value = "91.77"
$newSheet.Cells.Item(5, 5).Style = "Normal"
$newSheet.Cells.Item(5, 6).NumberFormat = "@"
$newSheet.Cells.Item(5, 6).Value = "8.21"
$newSheet.Cells.Item(5, 6).Style = "Normal"
$newSheet.Cells.Item(5, 7).NumberFormat = "@"
$newSheet.Cells.Item(5, 7).Value = "1.4269"
$newSheet.Cells.Item(5, 7).Style = "Normal"
$newSheet.Cells.Item(5, 8).Value = 4

# Row 6
$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).NumberFormat = "@"
$newSheet.Cells.Item(6, 2).Value = "009661"
$newSheet.Cells.Item(6, 2).Style = "Normal"
$newSheet.Cells.Item(6, 3).NumberFormat = "@"
$newSheet.Cells.Item(6, 3).Value = "平安研究睿选混合A"
$newSheet.Cells.Item(6, 3).Style = "Normal"
$newSheet.Cells.Item(6, 4).NumberFormat = "@"
$newSheet.Cells.Item(6, 4).Value = "22.07"
$newSheet.Cells.Item(6, 4).Style = "Normal"
$newSheet.Cells.Item(6, 5).NumberFormat = "@"
$newSheet.Cells.Item(6, 5).Value = "89.85"
$newSheet.Cells.Item(6, 5).Style = "Normal"
$newSheet.Cells.Item(6, 6).NumberFormat = "@"
$newSheet.Cells.Item(6, 6).Value = "4.95"
$newSheet.Cells.Item(6, 6).Style = "Normal"
$newSheet.Cells.Item(6, 7).NumberFormat = "@"
$newSheet.Cells.Item(6, 7).Value = "1.0925"
$newSheet.Cells.Item(6, 7).Style = "Normal"
$newSheet.Cells.Item(6, 8).Value = 8

# Row 7
$newSheet.Cells.Item(7, 1).Value = 5
$newSheet.Cells.Item(7, 2).NumberFormat = "@"
$newSheet.Cells.Item(7, 2).Value = "010242"
$newSheet.Cells.Item(7, 2).Style = "Normal"
$newSheet.Cells.Item(7, 3).NumberFormat = "@"
$newSheet.Cells.Item(7, 3).Value = "平安稳健增长混合A"
$newSheet.Cells.Item(7, 3).Style = "Normal"
$newSheet.Cells.Item(7, 4).NumberFormat = "@"
$newSheet.Cells.Item(7, 4).Value = "50.24"
$newSheet.Cells.Item(7, 4).Style = "Normal"
$newSheet.Cells.Item(7, 5).NumberFormat = "@"
$newSheet.Cells.Item(7, 5).Value = "34.03"
$newSheet.Cells.Item(7, 5).Style = "Normal"
$newSheet.Cells.Item(7, 6).NumberFormat = "@"
$newSheet.Cells.Item(7, 6).Value = "1.92"
$newSheet.Cells.Item(7, 6).Style = "Normal"
$newSheet.Cells.Item(7, 7).NumberFormat = "@"
$newSheet.Cells.Item(7, 7).Value = "0.9646"
$newSheet.Cells.Item(7, 7).Style = "Normal"
$newSheet.Cells.Item(7, 8).Value = 8

# Row 8
$newSheet.Cells.Item(8, 1).Value = 6
$newSheet.Cells.Item(8, 2).NumberFormat = "@"
$newSheet.Cells.Item(8, 2).Value = "009805"
$newSheet.Cells.Item(8, 2).Style = "Normal"
$newSheet.Cells.Item(8, 3).NumberFormat = "@"
$newSheet.Cells.Item(8, 3).Value = "国泰医药健康股票A"
$newSheet.Cells.Item(8, 3).Style = "Normal"
$newSheet.Cells.Item(8, 4).NumberFormat = "@"
$newSheet.Cells.Item(8, 4).Value = "12.14"
$newSheet.Cells.Item(8, 4).Style = "Normal"
$newSheet.Cells.Item(8, 5).NumberFormat = "@"
$newSheet.Cells.Item(8, 5).Value = "92.49"
$newSheet.Cells.Item(8, 5).Style = "Normal"
$newSheet.Cells.Item(8, 6).NumberFormat = "@"
$newSheet.Cells.Item(8, 6).Value = "5.53"
$newSheet.Cells.Item(8, 6).Style = "Normal"
$newSheet.Cells.Item(8, 7).NumberFormat = "@"
$newSheet.Cells.Item(8, 7).Value = "0.6713"
$newSheet.Cells.Item(8, 7).Style = "Normal"
$newSheet.Cells.Item(8, 8).Value = 9

# Row 9
$newSheet.Cells.Item(9, 1).Value = 7
$newSheet.Cells.Item(9, 2).NumberFormat = "@"
$newSheet.Cells.Item(9, 2).Value = "519908"
$newSheet.Cells.Item(9, 2).Style = "Normal"
$newSheet.Cells.Item(9, 3).NumberFormat = "@"
$newSheet.Cells.Item(9, 3).Value = "华夏兴华混合A"
$newSheet.Cells.Item(9, 3).Style = "Normal"
$newSheet.Cells.Item(9, 4).NumberFormat = "@"
$newSheet.Cells.Item(9, 4).Value = "9.39"
$newSheet.Cells.Item(9, 4).Style = "Normal"
$newSheet.Cells.Item(9, 5).NumberFormat = "@"
$newSheet.Cells.Item(9, 5).Value = "91.83"
$newSheet.Cells.Item(9, 5).Style = "Normal"
$newSheet.Cells.Item(9, 6).NumberFormat = "@"
$newSheet.Cells.Item(9, 6).Value = "6.12"
$newSheet.Cells.Item(9, 6).Style = "Normal"
$newSheet.Cells.Item(9, 7).NumberFormat = "@"
$newSheet.Cells.Item(9, 7).Value = "0.5747"
$newSheet.Cells.Item(9, 7).Style = "Normal"
$newSheet.Cells.Item(9, 8).Value = 2

# Row 10
$newSheet.Cells.Item(10, 1).Value = 8
$newSheet.Cells.Item(10, 2).NumberFormat = "@"
$newSheet.Cells.Item(10, 2).Value = "960004"
$newSheet.Cells.Item(10, 2).Style = "Normal"
$newSheet.Cells.Item(10, 3).NumberFormat = "@"
$newSheet.Cells.Item(10, 3).Value = "华夏兴华混合H"
$newSheet.Cells.Item(10, 3).Style = "Normal"
$newSheet.Cells.Item(10, 4).NumberFormat = "@"
$newSheet.Cells.Item(10, 4).Value = "9.39"
$newSheet.Cells.Item(10, 4).Style = "Normal"
$newSheet.Cells.Item(10, 5).NumberFormat = "@"
$newSheet.Cells.Item(10, 5).Value = "91.83"
$newSheet.Cells.Item(10, 5).Style = "Normal"
$newSheet.Cells.Item(10, 6).NumberFormat = "@"
$newSheet.Cells.Item(10, 6).Value = "6.12"
$newSheet.Cells.Item(10, 6).Style = "Normal"
$newSheet.Cells.Item(10, 7).NumberFormat = "@"
$newSheet.Cells.Item(10, 7).Value = "0.5747"
$newSheet.Cells.Item(10, 7).Style = "Normal"
$newSheet.Cells.Item(10, 8).Value = 2

# Row 11
$newSheet.Cells.Item(11, 1).Value = 9
$newSheet.Cells.Item(11, 2).NumberFormat = "@"
$newSheet.Cells.Item(11, 2).Value = "160215"
$newSheet.Cells.Item(11, 2).Style = "Normal"
$newSheet.Cells.Item(11, 3).NumberFormat = "@"
$newSheet.Cells.Item(11, 3).Value = "国泰价值经典灵活配置混合（LOF）"
$newSheet.Cells.Item(11, 3).Style = "Normal"
$newSheet.Cells.Item(11, 4).NumberFormat = "@"
$newSheet.Cells.Item(11, 4).Value = "6.36"
$newSheet.Cells.Item(11, 4).Style = "Normal"
$newSheet.Cells.Item(11, 5).NumberFormat = "@"
$newSheet.Cells.Item(11, 5).Value = "92.86"
$newSheet.Cells.Item(11, 5).Style = "Normal"
$newSheet.Cells.Item(11, 6).NumberFormat = "@"
$newSheet.Cells.Item(11, 6).Value = "8.20"
$newSheet.Cells.Item(11, 6).Style = "Normal"
$newSheet.Cells.Item(11, 7).NumberFormat = "@"
$newSheet.Cells.Item(11, 7).Value = "0.5215"
$newSheet.Cells.Item(11, 7).Style = "Normal"
$newSheet.Cells.Item(11, 8).Value = 4

# Row 12
$newSheet.Cells.Item(12, 1).Value = 10
$newSheet.Cells.Item(12, 2).NumberFormat = "@"
$newSheet.Cells.Item(12, 2).Value = "007082"
$newSheet.Cells.Item(12, 2).Style = "Normal"
$newSheet.Cells.Item(12, 3).NumberFormat = "@"
$newSheet.Cells.Item(12, 3).Value = "平安高端制造混合A"
$newSheet.Cells.Item(12, 3).Style = "Normal"
$newSheet.Cells.Item(12, 4).NumberFormat = "@"
$newSheet.Cells.Item(12, 4).Value = "10.45"
$newSheet.Cells.Item(12, 4).Style = "Normal"
$newSheet.Cells.Item(12, 5).NumberFormat = "@"
$newSheet.Cells.Item(12, 5).Value = "88.42"
$newSheet.Cells.Item(12, 5).Style = "Normal"
$newSheet.Cells.Item(12, 6).NumberFormat = "@"
$newSheet.Cells.Item(12, 6).Value = "4.89"
$newSheet.Cells.Item(12, 6).Style = "Normal"
$newSheet.Cells.Item(12, 7).NumberFormat = "@"
$newSheet.Cells.Item(12, 7).Value = "0.5110"
$newSheet.Cells.Item(12, 7).Style = "Normal"
$newSheet.Cells.Item(12, 8).Value = 7

# Row 13
$newSheet.Cells.Item(13, 1).Value = 11
$newSheet.Cells.Item(13, 2).NumberFormat = "@"
$newSheet.Cells.Item(13, 2).Value = "011392"
$newSheet.Cells.Item(13, 2).Style = "Normal"
$newSheet.Cells.Item(13, 3).NumberFormat = "@"
$newSheet.Cells.Item(13, 3).Value = "平安兴鑫回报一年定期开放混合"
$newSheet.Cells.Item(13, 3).Style = "Normal"
$newSheet.Cells.Item(13, 4).NumberFormat = "@"
$newSheet.Cells.Item(13, 4).Value = "7.07"
$newSheet.Cells.Item(13, 4).Style = "Normal"
$newSheet.Cells.Item(13, 5).NumberFormat = "@"
$newSheet.Cells.Item(13, 5).Value = "99.45"
$newSheet.Cells.Item(13, 5).Style = "Normal"
$newSheet.Cells.Item(13, 6).NumberFormat = "@"
$newSheet.Cells.Item(13, 6).Value = "6.21"
$newSheet.Cells.Item(13, 6).Style = "Normal"
$newSheet.Cells.Item(13, 7).NumberFormat = "@"
$newSheet.Cells.Item(13, 7).Value = "0.4390"
$newSheet.Cells.Item(13, 7).Style = "Normal"
$newSheet.Cells.Item(13, 8).Value = 6

# Row 14
$newSheet.Cells.Item(14, 1).Value = 12
$newSheet.Cells.Item(14, 2).NumberFormat = "@"
$newSheet.Cells.Item(14, 2).Value = "008370"
$newSheet.Cells.Item(14, 2).Style = "Normal"
$newSheet.Cells.Item(14, 3).NumberFormat = "@"
$newSheet.Cells.Item(14, 3).Value = "国泰研究精选两年持有期混合"
$newSheet.Cells.Item(14, 3).Style = "Normal"
$newSheet.Cells.Item(14, 4).NumberFormat = "@"
$newSheet.Cells.Item(14, 4).Value = "4.15"
$newSheet.Cells.Item(14, 4).Style = "Normal"
$newSheet.Cells.Item(14, 5).NumberFormat = "@"
$newSheet.Cells.Item(14, 5).Value = "92.87"
$newSheet.Cells.Item(14, 5).Style = "Normal"
$newSheet.Cells.Item(14, 6).NumberFormat = "@"
$newSheet.Cells.Item(14, 6).Value = "9.52"
$newSheet.Cells.Item(14, 6).Style = "Normal"
$newSheet.Cells.Item(14, 7).NumberFormat = "@"
$newSheet.Cells.Item(14, 7).Value = "0.3951"
$newSheet.Cells.Item(14, 7).Style = "Normal"
$newSheet.Cells.Item(14, 8).Value = 4

# Row 15
$newSheet.Cells.Item(15, 1).Value = 13
$newSheet.Cells.Item(15, 2).NumberFormat = "@"
$newSheet.Cells.Item(15, 2).Value = "010243"
$newSheet.Cells.Item(15, 2).Style = "Normal"
$newSheet.Cells.Item(15, 3).NumberFormat = "@"
$newSheet.Cells.Item(15, 3).Value = "平安稳健增长混合C"
$newSheet.Cells.Item(15, 3).Style = "Normal"
$newSheet.Cells.Item(15, 4).NumberFormat = "@"
$newSheet.Cells.Item(15, 4).Value = "17.89"
$newSheet.Cells.Item(15, 4).Style = "Normal"
$newSheet.Cells.Item(15, 5).NumberFormat = "@"
$newSheet.Cells.Item(15, 5).Value = "34.03"
$newSheet.Cells.Item(15, 5).Style = "Normal"
$newSheet.Cells.Item(15, 6).NumberFormat = "@"
$newSheet.Cells.Item(15, 6).Value = "1.92"
$newSheet.Cells.Item(15, 6).Style = "Normal"
$newSheet.Cells.Item(15, 7).NumberFormat = "@"
$newSheet.Cells.Item(15, 7).Value = "0.3435"
$newSheet.Cells.Item(15, 7).Style = "Normal"
$newSheet.Cells.Item(15, 8).Value = 8

# Row 16
$newSheet.Cells.Item(16, 1).Value = 14
$newSheet.Cells.Item(16, 2).NumberFormat = "@"
$newSheet.Cells.Item(16, 2).Value = "009804"
$newSheet.Cells.Item(16, 2).Style = "Normal"
$newSheet.Cells.Item(16, 3).NumberFormat = "@"
$newSheet.Cells.Item(16, 3).Value = "国泰研究优势混合"
$newSheet.Cells.Item(16, 3).Style = "Normal"
$newSheet.Cells.Item(16, 4).NumberFormat = "@"
$newSheet.Cells.Item(16, 4).Value = "4.14"
$newSheet.Cells.Item(16, 4).Style = "Normal"
$newSheet.Cells.Item(16, 5).NumberFormat = "@"
$newSheet.Cells.Item(16, 5).Value = "90.91"
$newSheet.Cells.Item(16, 5).Style = "Normal"
$newSheet.Cells.Item(16, 6).NumberFormat = "@"
$newSheet.Cells.Item(16, 6).Value = "8.11"
$newSheet.Cells.Item(16, 6).Style = "Normal"
$newSheet.Cells.Item(16, 7).NumberFormat = "@"
$newSheet.Cells.Item(16, 7).Value = "0.3358"
$newSheet.Cells.Item(16, 7).Style = "Normal"
$newSheet.Cells.Item(16, 8).Value = 6

# Row 17
$newSheet.Cells.Item(17, 1).Value = 15
$newSheet.Cells.Item(17, 2).NumberFormat = "@"
$newSheet.Cells.Item(17, 2).Value = "160212"
$newSheet.Cells.Item(17, 2).Style = "Normal"
$newSheet.Cells.Item(17, 3).NumberFormat = "@"
$newSheet.Cells.Item(17, 3).Value = "国泰估值优势混合 (LOF)"
$newSheet.Cells.Item(17, 3).Style = "Normal"
$newSheet.Cells.Item(17, 4).NumberFormat = "@"
$newSheet.Cells.Item(17, 4).Value = "8.98"
$newSheet.Cells.Item(17, 4).Style = "Normal"
$newSheet.Cells.Item(17, 5).NumberFormat = "@"
$newSheet.Cells.Item(17, 5).Value = "62.69"
$newSheet.Cells.Item(17, 5).Style = "Normal"
$newSheet.Cells.Item(17, 6).NumberFormat = "@"
$newSheet.Cells.Item(17, 6).Value = "3.50"
$newSheet.Cells.Item(17, 6).Style = "Normal"
$newSheet.Cells.Item(17, 7).NumberFormat = "@"
$newSheet.Cells.Item(17, 7).Value = "0.3143"
$newSheet.Cells.Item(17, 7).Style = "Normal"
$newSheet.Cells.Item(17, 8).Value = 8

# Row 18
$newSheet.Cells.Item(18, 1).Value = 16
$newSheet.Cells.Item(18, 2).NumberFormat = "@"
$newSheet.Cells.Item(18, 2).Value = "001297"
$newSheet.Cells.Item(18, 2).Style = "Normal"
$newSheet.Cells.Item(18, 3).NumberFormat = "@"
$newSheet.Cells.Item(18, 3).Value = "平安智慧中国灵活配置混合"
$newSheet.Cells.Item(18, 3).Style = "Normal"
$newSheet.Cells.Item(18, 4).NumberFormat = "@"
$newSheet.Cells.Item(18, 4).Value = "5.57"
$newSheet.Cells.Item(18, 4).Style = "Normal"
$newSheet.Cells.Item(18, 5).NumberFormat = "@"
$newSheet.Cells.Item(18, 5).Value = "84.03"
$newSheet.Cells.Item(18, 5).Style = "Normal"
$newSheet.Cells.Item(18, 6).NumberFormat = "@"
$newSheet.Cells.Item(18, 6).Value = "4.99"
$newSheet.Cells.Item(18, 6).Style = "Normal"
$newSheet.Cells.Item(18, 7).NumberFormat = "@"
$newSheet.Cells.Item(18, 7).Value = "0.2779"
$newSheet.Cells.Item(18, 7).Style = "Normal"
$newSheet.Cells.Item(18, 8).Value = 7

# Row 19
$newSheet.Cells.Item(19, 1).Value = 17
$newSheet.Cells.Item(19, 2).NumberFormat = "@"
$newSheet.Cells.Item(19, 2).Value = "011321"
$newSheet.Cells.Item(19, 2).Style = "Normal"
$newSheet.Cells.Item(19, 3).NumberFormat = "@"
$newSheet.Cells.Item(19, 3).Value = "国泰大健康股票C"
$newSheet.Cells.Item(19, 3).Style = "Normal"
$newSheet.Cells.Item(19, 4).NumberFormat = "@"
$newSheet.Cells.Item(19, 4).Value = "3.47"
$newSheet.Cells.Item(19, 4).Style = "Normal"
$newSheet.Cells.Item(19, 5).NumberFormat = "@"
$newSheet.Cells.Item(19, 5).Value = "90.83"
$newSheet.Cells.Item(19, 5).Style = "Normal"
$newSheet.Cells.Item(19, 6).NumberFormat = "@"
$newSheet.Cells.Item(19, 6).Value = "7.96"
$newSheet.Cells.Item(19, 6).Style = "Normal"
$newSheet.Cells.Item(19, 7).NumberFormat = "@"
$newSheet.Cells.Item(19, 7).Value = "0.2762"
$newSheet.Cells.Item(19, 7).Style = "Normal"
$newSheet.Cells.Item(19, 8).Value = 3

# Row 20
$newSheet.Cells.Item(20, 1).Value = 18
$newSheet.Cells.Item(20, 2).NumberFormat = "@"
$newSheet.Cells.Item(20, 2).Value = "013023"
$newSheet.Cells.Item(20, 2).Style = "Normal"
$newSheet.Cells.Item(20, 3).NumberFormat = "@"
$newSheet.Cells.Item(20, 3).Value = "平安均衡优选1年持有混合A"
$newSheet.Cells.Item(20, 3).Style = "Normal"
$newSheet.Cells.Item(20, 4).NumberFormat = "@"
$newSheet.Cells.Item(20, 4).Value = "4.16"
$newSheet.Cells.Item(20, 4).Style = "Normal"
$newSheet.Cells.Item(20, 5).NumberFormat = "@"
$newSheet.Cells.Item(20, 5).Value = "90.64"
$newSheet.Cells.Item(20, 5).Style = "Normal"
$newSheet.Cells.Item(20, 6).NumberFormat = "@"
$newSheet.Cells.Item(20, 6).Value = "4.95"
$newSheet.Cells.Item(20, 6).Style = "Normal"
$newSheet.Cells.Item(20, 7).NumberFormat = "@"
$newSheet.Cells.Item(20, 7).Value = "0.2059"
$newSheet.Cells.Item(20, 7).Style = "Normal"
$newSheet.Cells.Item(20, 8).Value = 7

# Row 21
$newSheet.Cells.Item(21, 1).Value = 19
$newSheet.Cells.Item(21, 2).NumberFormat = "@"
$newSheet.Cells.Item(21, 2).Value = "009662"
$newSheet.Cells.Item(21, 2).Style = "Normal"
$newSheet.Cells.Item(21, 3).NumberFormat = "@"
$newSheet.Cells.Item(21, 3).Value = "平安研究睿选混合C"
$newSheet.Cells.Item(21, 3).Style = "Normal"
$newSheet.Cells.Item(21, 4).NumberFormat = "@"
$newSheet.Cells.Item(21, 4).Value = "4.01"
$newSheet.Cells.Item(21, 4).Style = "Normal"
$newSheet.Cells.Item(21, 5).NumberFormat = "@"
$newSheet.Cells.Item(21, 5).Value = "89.85"
$newSheet.Cells.Item(21, 5).Style = "Normal"
$newSheet.Cells.Item(21, 6).NumberFormat = "@"
$newSheet.Cells.Item(21, 6).Value = "4.95"
$newSheet.Cells.Item(21, 6).Style = "Normal"
$newSheet.Cells.Item(21, 7).NumberFormat = "@"
$newSheet.Cells.Item(21, 7).Value = "0.1985"
$newSheet.Cells.Item(21, 7).Style = "Normal"
$newSheet.Cells.Item(21, 8).Value = 8

# Row 22
$newSheet.Cells.Item(22, 1).Value = 20
$newSheet.Cells.Item(22, 2).NumberFormat = "@"
$newSheet.Cells.Item(22, 2).Value = "007083"
$newSheet.Cells.Item(22, 2).Style = "Normal"
$newSheet.Cells.Item(22, 3).NumberFormat = "@"
$newSheet.Cells.Item(22, 3).Value = "平安高端制造混合C"
$newSheet.Cells.Item(22, 3).Style = "Normal"
$newSheet.Cells.Item(22, 4).NumberFormat = "@"
$newSheet.Cells.Item(22, 4).Value = "4.03"
$newSheet.Cells.Item(22, 4).Style = "Normal"
$newSheet.Cells.Item(22, 5).NumberFormat = "@"
$newSheet.Cells.Item(22, 5).Value = "88.42"
$newSheet.Cells.Item(22, 5).Style = "Normal"
$newSheet.Cells.Item(22, 6).NumberFormat = "@"
$newSheet.Cells.Item(22, 6).Value = "4.89"
$newSheet.Cells.Item(22, 6).Style = "Normal"
$newSheet.Cells.Item(22, 7).NumberFormat = "@"
$newSheet.Cells.Item(22, 7).Value = "0.1971"
$newSheet.Cells.Item(22, 7).Style = "Normal"
$newSheet.Cells.Item(22, 8).Value = 7

# Row 23
$newSheet.Cells.Item(23, 1).Value = 21
$newSheet.Cells.Item(23, 2).NumberFormat = "@"
$newSheet.Cells.Item(23, 2).Value = "001366"
$newSheet.Cells.Item(23, 2).Style = "Normal"
$newSheet.Cells.Item(23, 3).NumberFormat = "@"
$newSheet.Cells.Item(23, 3).Value = "金鹰产业整合灵活配置混合"
$newSheet.Cells.Item(23, 3).Style = "Normal"
$newSheet.Cells.Item(23, 4).NumberFormat = "@"
$newSheet.Cells.Item(23, 4).Value = "5.27"
$newSheet.Cells.Item(23, 4).Style = "Normal"
$newSheet.Cells.Item(23, 5).NumberFormat = "@"
$newSheet.Cells.Item(23, 5).Value = "91.56"
$newSheet.Cells.Item(23, 5).Style = "Normal"
$newSheet.Cells.Item(23, 6).NumberFormat = "@"
$newSheet.Cells.Item(23, 6).Value = "3.28"
$newSheet.Cells.Item(23, 6).Style = "Normal"
$newSheet.Cells.Item(23, 7).NumberFormat = "@"
$newSheet.Cells.Item(23, 7).Value = "0.1729"
$newSheet.Cells.Item(23, 7).Style = "Normal"
$newSheet.Cells.Item(23, 8).Value = 7

# Row 24
$newSheet.Cells.Item(24, 1).Value = 22
$newSheet.Cells.Item(24, 2).NumberFormat = "@"
$newSheet.Cells.Item(24, 2).Value = "008066"
$newSheet.Cells.Item(24, 2).Style = "Normal"
$newSheet.Cells.Item(24, 3).NumberFormat = "@"
$newSheet.Cells.Item(24, 3).Value = "汇添富中盘积极成长混合C"
$newSheet.Cells.Item(24, 3).Style = "Normal"
$newSheet.Cells.Item(24, 4).NumberFormat = "@"
$newSheet.Cells.Item(24, 4).Value = "4.31"
$newSheet.Cells.Item(24, 4).Style = "Normal"
$newSheet.Cells.Item(24, 5).NumberFormat = "@"
$newSheet.Cells.Item(24, 5).Value = "73.48"
$newSheet.Cells.Item(24, 5).Style = "Normal"
$newSheet.Cells.Item(24, 6).NumberFormat = "@"
$newSheet.Cells.Item(24, 6).Value = "3.68"
$newSheet.Cells.Item(24, 6).Style = "Normal"
$newSheet.Cells.Item(24, 7).NumberFormat = "@"
$newSheet.Cells.Item(24, 7).Value = "0.1586"
$newSheet.Cells.Item(24, 7).Style = "Normal"
$newSheet.Cells.Item(24, 8).Value = 4

# Row 25
$newSheet.Cells.Item(25, 1).Value = 23
$newSheet.Cells.Item(25, 2).NumberFormat = "@"
$newSheet.Cells.Item(25, 2).Value = "011326"
$newSheet.Cells.Item(25, 2).Style = "Normal"
$newSheet.Cells.Item(25, 3).NumberFormat = "@"
$newSheet.Cells.Item(25, 3).Value = "国泰医药健康股票C"
$newSheet.Cells.Item(25, 3).Style = "Normal"
$newSheet.Cells.Item(25, 4).NumberFormat = "@"
$newSheet.Cells.Item(25, 4).Value = "1.09"
$newSheet.Cells.Item(25, 4).Style = "Normal"
$newSheet.Cells.Item(25, 5).NumberFormat = "@"
$newSheet.Cells.Item(25, 5).Value = "92.49"
$newSheet.Cells.Item(25, 5).Style = "Normal"
$newSheet.Cells.Item(25, 6).NumberFormat = "@"
$newSheet.Cells.Item(25, 6).Value = "5.53"
$newSheet.Cells.Item(25, 6).Style = "Normal"
$newSheet.Cells.Item(25, 7).NumberFormat = "@"
$newSheet.Cells.Item(25, 7).Value = "0.0603"
$newSheet.Cells.Item(25, 7).Style = "Normal"
$newSheet.Cells.Item(25, 8).Value = 9

# Row 26
$newSheet.Cells.Item(26, 1).Value = 24
$newSheet.Cells.Item(26, 2).NumberFormat = "@"
$newSheet.Cells.Item(26, 2).Value = "004266"
$newSheet.Cells.Item(26, 2).Style = "Normal"
$newSheet.Cells.Item(26, 3).NumberFormat = "@"
$newSheet.Cells.Item(26, 3).Value = "招商沪港深科技创新主题精选灵活配置混合A"
$newSheet.Cells.Item(26, 3).Style = "Normal"
$newSheet.Cells.Item(26, 4).NumberFormat = "@"
$newSheet.Cells.Item(26, 4).Value = "1.29"
$newSheet.Cells.Item(26, 4).Style = "Normal"
$newSheet.Cells.Item(26, 5).NumberFormat = "@"
$newSheet.Cells.Item(26, 5).Value = "88.85"
$newSheet.Cells.Item(26, 5).Style = "Normal"
$newSheet.Cells.Item(26, 6).NumberFormat = "@"
$newSheet.Cells.Item(26, 6).Value = "3.69"
$newSheet.Cells.Item(26, 6).Style = "Normal"
$newSheet.Cells.Item(26, 7).NumberFormat = "@"
$newSheet.Cells.Item(26, 7).Value = "0.0476"
$newSheet.Cells.Item(26, 7).Style = "Normal"
$newSheet.Cells.Item(26, 8).Value = 5

# Row 27
$newSheet.Cells.Item(27, 1).Value = 25
$newSheet.Cells.Item(27, 2).NumberFormat = "@"
$newSheet.Cells.Item(27, 2).Value = "210006"
$newSheet.Cells.Item(27, 2).Style = "Normal"
$newSheet.Cells.Item(27, 3).NumberFormat = "@"
$newSheet.Cells.Item(27, 3).Value = "金鹰元禧混合A"
$newSheet.Cells.Item(27, 3).Style = "Normal"
$newSheet.Cells.Item(27, 4).NumberFormat = "@"
$newSheet.Cells.Item(27, 4).Value = "7.07"
$newSheet.Cells.Item(27, 4).Style = "Normal"
$newSheet.Cells.Item(27, 5).NumberFormat = "@"
$newSheet.Cells.Item(27, 5).Value = "22.58"
$newSheet.Cells.Item(27, 5).Style = "Normal"
$newSheet.Cells.Item(27, 6).NumberFormat = "@"
$newSheet.Cells.Item(27, 6).Value = "0.52"
$newSheet.Cells.Item(27, 6).Style = "Normal"
$newSheet.Cells.Item(27, 7).NumberFormat = "@"
$newSheet.Cells.Item(27, 7).Value = "0.0368"
$newSheet.Cells.Item(27, 7).Style = "Normal"
$newSheet.Cells.Item(27, 8).Value = 5

# Row 28
$newSheet.Cells.Item(28, 1).Value = 26
$newSheet.Cells.Item(28, 2).NumberFormat = "@"
$newSheet.Cells.Item(28, 2).Value = "000110"
$newSheet.Cells.Item(28, 2).Style = "Normal"
$newSheet.Cells.Item(28, 3).NumberFormat = "@"
$newSheet.Cells.Item(28, 3).Value = "金鹰元安混合A"
$newSheet.Cells.Item(28, 3).Style = "Normal"
$newSheet.Cells.Item(28, 4).NumberFormat = "@"
$newSheet.Cells.Item(28, 4).Value = "7.13"
$newSheet.Cells.Item(28, 4).Style = "Normal"
$newSheet.Cells.Item(28, 5).NumberFormat = "@"
$newSheet.Cells.Item(28, 5).Value = "22.33"
$newSheet.Cells.Item(28, 5).Style = "Normal"
$newSheet.Cells.Item(28, 6).NumberFormat = "@"
$newSheet.Cells.Item(28, 6).Value = "0.49"
$newSheet.Cells.Item(28, 6).Style = "Normal"
$newSheet.Cells.Item(28, 7).NumberFormat = "@"
$newSheet.Cells.Item(28, 7).Value = "0.0349"
$newSheet.Cells.Item(28, 7).Style = "Normal"
$newSheet.Cells.Item(28, 8).Value = 4

# Row 29
$newSheet.Cells.Item(29, 1).Value = 27
$newSheet.Cells.Item(29, 2).NumberFormat = "@"
$newSheet.Cells.Item(29, 2).Value = "210010"
$newSheet.Cells.Item(29, 2).Style = "Normal"
$newSheet.Cells.Item(29, 3).NumberFormat = "@"
$newSheet.Cells.Item(29, 3).Value = "金鹰灵活配置混合A"
$newSheet.Cells.Item(29, 3).Style = "Normal"
$newSheet.Cells.Item(29, 4).NumberFormat = "@"
$newSheet.Cells.Item(29, 4).Value = "6.15"
$newSheet.Cells.Item(29, 4).Style = "Normal"
$newSheet.Cells.Item(29, 5).NumberFormat = "@"
$newSheet.Cells.Item(29, 5).Value = "23.34"
$newSheet.Cells.Item(29, 5).Style = "Normal"
$newSheet.Cells.Item(29, 6).NumberFormat = "@"
$newSheet.Cells.Item(29, 6).Value = "0.55"
$newSheet.Cells.Item(29, 6).Style = "Normal"
$newSheet.Cells.Item(29, 7).NumberFormat = "@"
$newSheet.Cells.Item(29, 7).Value = "0.0338"
$newSheet.Cells.Item(29, 7).Style = "Normal"
$newSheet.Cells.Item(29, 8).Value = 8

# Row 30
$newSheet.Cells.Item(30, 1).Value = 28
$newSheet.Cells.Item(30, 2).NumberFormat = "@"
$newSheet.Cells.Item(30, 2).Value = "210011"
$newSheet.Cells.Item(30, 2).Style = "Normal"
$newSheet.Cells.Item(30, 3).NumberFormat = "@"
$newSheet.Cells.Item(30, 3).Value = "金鹰灵活配置混合C"
$newSheet.Cells.Item(30, 3).Style = "Normal"
$newSheet.Cells.Item(30, 4).NumberFormat = "@"
$newSheet.Cells.Item(30, 4).Value = "3.61"
$newSheet.Cells.Item(30, 4).Style = "Normal"
$newSheet.Cells.Item(30, 5).NumberFormat = "@"
$newSheet.Cells.Item(30, 5).Value = "23.34"
$newSheet.Cells.Item(30, 5).Style = "Normal"
$newSheet.Cells.Item(30, 6).NumberFormat = "@"
$newSheet.Cells.Item(30, 6).Value = "0.55"
$newSheet.Cells.Item(30, 6).Style = "Normal"
$newSheet.Cells.Item(30, 7).NumberFormat = "@"
$newSheet.Cells.Item(30, 7).Value = "0.0199"
$newSheet.Cells.Item(30, 7).Style = "Normal"
$newSheet.Cells.Item(30, 8).Value = 8

# Row 31
$newSheet.Cells.Item(31, 1).Value = 29
$newSheet.Cells.Item(31, 2).NumberFormat = "@"
$newSheet.Cells.Item(31, 2).Value = "002425"
$newSheet.Cells.Item(31, 2).Style = "Normal"
$newSheet.Cells.Item(31, 3).NumberFormat = "@"
$newSheet.Cells.Item(31, 3).Value = "金鹰元禧混合C"
$newSheet.Cells.Item(31, 3).Style = "Normal"
$newSheet.Cells.Item(31, 4).NumberFormat = "@"
$newSheet.Cells.Item(31, 4).Value = "3.19"
$newSheet.Cells.Item(31, 4).Style = "Normal"
$newSheet.Cells.Item(31, 5).NumberFormat = "@"
$newSheet.Cells.Item(31, 5).Value = "22.58"
$newSheet.Cells.Item(31, 5).Style = "Normal"
$newSheet.Cells.Item(31, 6).NumberFormat = "@"
$newSheet.Cells.Item(31, 6).Value = "0.52"
$newSheet.Cells.Item(31, 6).Style = "Normal"
$newSheet.Cells.Item(31, 7).NumberFormat = "@"
$newSheet.Cells.Item(31, 7).Value = "0.0166"
$newSheet.Cells.Item(31, 7).Style = "Normal"
$newSheet.Cells.Item(31, 8).Value = 5

# Row 32
$newSheet.Cells.Item(32, 1).Value = 30
$newSheet.Cells.Item(32, 2).NumberFormat = "@"
$newSheet.Cells.Item(32, 2).Value = "013024"
$newSheet.Cells.Item(32, 2).Style = "Normal"
$newSheet.Cells.Item(32, 3).NumberFormat = "@"
$newSheet.Cells.Item(32, 3).Value = "平安均衡优选1年持有混合C"
$newSheet.Cells.Item(32, 3).Style = "Normal"
$newSheet.Cells.Item(32, 4).NumberFormat = "@"
$newSheet.Cells.Item(32, 4).Value = "0.25"
$newSheet.Cells.Item(32, 4).Style = "Normal"
$newSheet.Cells.Item(32, 5).NumberFormat = "@"
$newSheet.Cells.Item(32, 5).Value = "90.64"
$newSheet.Cells.Item(32, 5).Style = "Normal"
$newSheet.Cells.Item(32, 6).NumberFormat = "@"
$newSheet.Cells.Item(32, 6).Value = "4.95"
$newSheet.Cells.Item(32, 6).Style = "Normal"
$newSheet.Cells.Item(32, 7).NumberFormat = "@"
$newSheet.Cells.Item(32, 7).Value = "0.0124"
$newSheet.Cells.Item(32, 7).Style = "Normal"
$newSheet.Cells.Item(32, 8).Value = 7

# Row 33
$newSheet.Cells.Item(33, 1).Value = 31
$newSheet.Cells.Item(33, 2).NumberFormat = "@"
$newSheet.Cells.Item(33, 2).Value = "002513"
$newSheet.Cells.Item(33, 2).Style = "Normal"
$newSheet.Cells.Item(33, 3).NumberFormat = "@"
$newSheet.Cells.Item(33, 3).Value = "金鹰元安混合C"
$newSheet.Cells.Item(33, 3).Style = "Normal"
$newSheet.Cells.Item(33, 4).NumberFormat = "@"
$newSheet.Cells.Item(33, 4).Value = "2.46"
$newSheet.Cells.Item(33, 4).Style = "Normal"
$newSheet.Cells.Item(33, 5).NumberFormat = "@"
$newSheet.Cells.Item(33, 5).Value = "22.33"
$newSheet.Cells.Item(33, 5).Style = "Normal"
$newSheet.Cells.Item(33, 6).NumberFormat = "@"
$newSheet.Cells.Item(33, 6).Value = "0.49"
$newSheet.Cells.Item(33, 6).Style = "Normal"
$newSheet.Cells.Item(33, 7).NumberFormat = "@"
$newSheet.Cells.Item(33, 7).Value = "0.0121"
$newSheet.Cells.Item(33, 7).Style = "Normal"
$newSheet.Cells.Item(33, 8).Value = 4

# Row 34
$newSheet.Cells.Item(34, 1).Value = 32
$newSheet.Cells.Item(34, 2).NumberFormat = "@"
$newSheet.Cells.Item(34, 2).Value = "010754"
$newSheet.Cells.Item(34, 2).Style = "Normal"
$newSheet.Cells.Item(34, 3).NumberFormat = "@"
$newSheet.Cells.Item(34, 3).Value = "招商沪港深科技创新主题精选灵活配置混合C"
$newSheet.Cells.Item(34, 3).Style = "Normal"
$newSheet.Cells.Item(34, 4).NumberFormat = "@"
$newSheet.Cells.Item(34, 4).Value = "0.28"
$newSheet.Cells.Item(34, 4).Style = "Normal"
$newSheet.Cells.Item(34, 5).NumberFormat = "@"
$newSheet.Cells.Item(34, 5).Value = "88.85"
$newSheet.Cells.Item(34, 5).Style = "Normal"
$newSheet.Cells.Item(34, 6).NumberFormat = "@"
$newSheet.Cells.Item(34, 6).Value = "3.69"
$newSheet.Cells.Item(34, 6).Style = "Normal"
$newSheet.Cells.Item(34, 7).NumberFormat = "@"
$newSheet.Cells.Item(34, 7).Value = "0.0103"
$newSheet.Cells.Item(34, 7).Style = "Normal"
$newSheet.Cells.Item(34, 8).Value = 5

# 2. Update the "总计" (total) summary sheet: insert a new row for 2022-Q1 at the top of the data.
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 33
$totalSheet.Cells.Item(2, 4).Value = 16.41

# Fix up formatting on the inserted row: column A keeps the bold/bordered "index" style,
# columns B:D should have no special style (matching the other data rows).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

# Renumber the running index in column A for the rows pushed down by the insert.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

# Restore the original active sheet/selection so the workbook view state is unchanged.
$wb.Worksheets.Item("2020-Q4").Select()
$wb.Worksheets.Item("2020-Q4").Range("A1").Select()

